$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.276.40"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.241.62"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.32"
$ws.Range("E5").Value = "  -2.63%  "

$ws.Range("E6").Value = "  -4.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").Value = "  -0.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("E9").Value = "  -1.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.56"
$ws.Range("E10").Value = "  -4.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0811"
$ws.Range("E11").Value = "  -1.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.17"
$ws.Range("E12").Value = "  -2.22%  "

$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.335.86"
$ws.Range("E14").Value = "  +3.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.582.45"
$ws.Range("E15").Value = "  +0.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.829"
$ws.Range("E16").Value = "  -1.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.50"
$ws.Range("E17").Value = "  -3.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.022.29"
$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0963"
$ws.Range("E19").Value = "  -1.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.39"
$ws.Range("E20").Value = "  +0.95%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.18"
$ws.Range("E21").Value = "  -7.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.51"
$ws.Range("E22").Value = "  -0.16%  "

$ws.Range("E23").Value = "  +0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.94"
$ws.Range("E24").Value = "  -1.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.00"
$ws.Range("E25").Value = "  -1.01%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.83"
$ws.Range("E27").Value = "  +6.92%  "

$ws.Range("E28").Value = "  -1.87%  "

$ws.Range("E29").Value = "  +4.16%  "

$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("E31").Value = "  -1.94%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.42"
$ws.Range("E32").Value = "  -1.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0793"
$ws.Range("E33").Value = "  -5.31%  "

$ws.Range("E34").Value = "  -1.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.13"
$ws.Range("E35").Value = "  -4.30%  "

$ws.Range("E36").Value = "  +1.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.108"
$ws.Range("E37").Value = "  +0.22%  "

$ws.Range("E38").Value = "  -7.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.49"
$ws.Range("E39").Value = "  -0.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.80"
$ws.Range("E40").Value = "  -4.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.27"
$ws.Range("E41").Value = "  -7.13%  "

$ws.Range("E42").Value = "  -2.74%  "

$ws.Range("E43").Value = "  +0.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.746.92"
$ws.Range("E44").Value = "  +2.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "82.93"
$ws.Range("E45").Value = "  +0.57%  "

$ws.Range("E46").Value = "  -1.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.61"
$ws.Range("E47").Value = "  -1.90%  "

$ws.Range("E48").Value = "  -4.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.60"
$ws.Range("E49").Value = "  +0.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.07"
$ws.Range("E50").Value = "  -0.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.69"
$ws.Range("E51").Value = "  -2.77%  "
